# Daily attendance processing - swap the ordering of names in the
# "Recorded By" column (column G) on the "Session Analysis Results" sheet.
#
# Two specific value patterns are reordered:
#   "dnasr281@gmail.com, System"            -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, system, System"   -> "backup@backdoor.com, System, system"
# All other values in column G are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = 7
    $val = $cell.Value()

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
}
